# TC001.xlsx edit — 19th September 2017, 06:33 PM
#
# Changes made on the TC001 sheet:
#  1. C6 ("Data" column for the "enterByXpath" row on the password field)
#     is updated from the old dummy password value to "rohith".
#  2. The active cell / selection on the sheet moves from C8 to B20.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC001")

$ws.Range("C6").Value = "rohith"

$ws.Range("B20").Select()
